$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows from 28 data rows (A1:C29) to 33 data rows (A1:C34).
# Give the 5 new rows the same formatting as the existing data rows by
# copying the format of the last existing data row down onto them.
$ws.Range("A29:C29").Copy()
$ws.Range("A30:C34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Keep phone numbers, DDD codes and the date strings as literal text
# (not auto-converted to numbers / date serials) for all data rows.
$ws.Range("A2:C34").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2,1).Value = '+5514997883211'
$ws.Cells.Item(2,2).Value = '14'
$ws.Cells.Item(2,3).Value = '2024-10-31'

# Row 3
$ws.Cells.Item(3,1).Value = '+5511920075911'
$ws.Cells.Item(3,2).Value = '11'
$ws.Cells.Item(3,3).Value = '2024-10-25'

# Row 4
$ws.Cells.Item(4,1).Value = '+5514996538302'
$ws.Cells.Item(4,2).Value = '14'
$ws.Cells.Item(4,3).Value = '2024-10-25'

# Row 5
$ws.Cells.Item(5,1).Value = '+556184978538'
$ws.Cells.Item(5,2).Value = '61'
$ws.Cells.Item(5,3).Value = '2024-10-22'

# Row 6
$ws.Cells.Item(6,1).Value = '+5511964652979'
$ws.Cells.Item(6,2).Value = '11'
$ws.Cells.Item(6,3).Value = '2024-10-21'

# Row 7
$ws.Cells.Item(7,1).Value = '+5521964219027'
$ws.Cells.Item(7,2).Value = '21'
$ws.Cells.Item(7,3).Value = '2024-10-18'

# Row 8
$ws.Cells.Item(8,1).Value = '+5513997050892'
$ws.Cells.Item(8,2).Value = '13'
$ws.Cells.Item(8,3).Value = '2024-10-18'

# Row 9
$ws.Cells.Item(9,1).Value = '+5519997201600'
$ws.Cells.Item(9,2).Value = '19'
$ws.Cells.Item(9,3).Value = '2024-10-18'

# Row 10
$ws.Cells.Item(10,1).Value = '+5517991290893'
$ws.Cells.Item(10,2).Value = '17'
$ws.Cells.Item(10,3).Value = '2024-10-17'

# Row 11
$ws.Cells.Item(11,1).Value = '+5511975292030'
$ws.Cells.Item(11,2).Value = '11'
$ws.Cells.Item(11,3).Value = '2024-10-14'

# Row 12
$ws.Cells.Item(12,1).Value = '+555199199744'
$ws.Cells.Item(12,2).Value = '51'
$ws.Cells.Item(12,3).Value = '2024-10-10'

# Row 13
$ws.Cells.Item(13,1).Value = '+5515996313442'
$ws.Cells.Item(13,2).Value = '15'
$ws.Cells.Item(13,3).Value = '2024-10-09'

# Row 14
$ws.Cells.Item(14,1).Value = '+5522981222545'
$ws.Cells.Item(14,2).Value = '22'
$ws.Cells.Item(14,3).Value = '2024-10-05'

# Row 15
$ws.Cells.Item(15,1).Value = '+553791176954'
$ws.Cells.Item(15,2).Value = '37'
$ws.Cells.Item(15,3).Value = '2024-10-04'

# Row 16
$ws.Cells.Item(16,1).Value = '+5516993020307'
$ws.Cells.Item(16,2).Value = '16'
$ws.Cells.Item(16,3).Value = '2024-10-04'

# Row 17
$ws.Cells.Item(17,1).Value = '+555491557534'
$ws.Cells.Item(17,2).Value = '54'
$ws.Cells.Item(17,3).Value = '2024-10-04'

# Row 18
$ws.Cells.Item(18,1).Value = '+556181971614'
$ws.Cells.Item(18,2).Value = '61'
$ws.Cells.Item(18,3).Value = '2024-10-01'

# Row 19
$ws.Cells.Item(19,1).Value = '+5513988453610'
$ws.Cells.Item(19,2).Value = '13'
$ws.Cells.Item(19,3).Value = '2024-09-25'

# Row 20
$ws.Cells.Item(20,1).Value = '+5511966134418'
$ws.Cells.Item(20,2).Value = '11'
$ws.Cells.Item(20,3).Value = '2024-09-24'

# Row 21
$ws.Cells.Item(21,1).Value = '+5511967085107'
$ws.Cells.Item(21,2).Value = '11'
$ws.Cells.Item(21,3).Value = '2024-09-20'

# Row 22
$ws.Cells.Item(22,1).Value = '+5511977696904'
$ws.Cells.Item(22,2).Value = '11'
$ws.Cells.Item(22,3).Value = '2024-09-20'

# Row 23
$ws.Cells.Item(23,1).Value = '+556192771804'
$ws.Cells.Item(23,2).Value = '61'
$ws.Cells.Item(23,3).Value = '2024-09-18'

# Row 24
$ws.Cells.Item(24,1).Value = '+556198454144'
$ws.Cells.Item(24,2).Value = '61'
$ws.Cells.Item(24,3).Value = '2024-09-18'

# Row 25
$ws.Cells.Item(25,1).Value = '+5511952381413'
$ws.Cells.Item(25,2).Value = '11'
$ws.Cells.Item(25,3).Value = '2024-09-18'

# Row 26
$ws.Cells.Item(26,1).Value = '+555199100909'
$ws.Cells.Item(26,2).Value = '51'
$ws.Cells.Item(26,3).Value = '2024-09-13'

# Row 27
$ws.Cells.Item(27,1).Value = '+5521965489343'
$ws.Cells.Item(27,2).Value = '21'
$ws.Cells.Item(27,3).Value = '2024-09-11'

# Row 28
$ws.Cells.Item(28,1).Value = '+5511947261969'
$ws.Cells.Item(28,2).Value = '11'
$ws.Cells.Item(28,3).Value = '2024-09-10'

# Row 29
$ws.Cells.Item(29,1).Value = '+5521997432262'
$ws.Cells.Item(29,2).Value = '21'
$ws.Cells.Item(29,3).Value = '2024-09-10'

# Row 30
$ws.Cells.Item(30,1).Value = '+5521985109311'
$ws.Cells.Item(30,2).Value = '21'
$ws.Cells.Item(30,3).Value = '2024-09-09'

# Row 31
$ws.Cells.Item(31,1).Value = '+553291004823'
$ws.Cells.Item(31,2).Value = '32'
$ws.Cells.Item(31,3).Value = '2024-08-26'

# Row 32
$ws.Cells.Item(32,1).Value = '+34603138909'
$ws.Cells.Item(32,2).Value = ''
$ws.Cells.Item(32,3).Value = '2024-08-19'

# Row 33
$ws.Cells.Item(33,1).Value = '+5511967859426'
$ws.Cells.Item(33,2).Value = '11'
$ws.Cells.Item(33,3).Value = '2024-07-28'

# Row 34
$ws.Cells.Item(34,1).Value = '+556298529715'
$ws.Cells.Item(34,2).Value = '62'
$ws.Cells.Item(34,3).Value = '2024-07-09'
